$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# The "review_info" sheet erroneously carried one data row (hotel data) - remove it,
# leaving only the header row.
$reviewSheet.Rows(2).Delete()

# The "hotel_info" sheet is missing a "State" column between "Hotel_Name" and "City".
# Insert it and populate the new column.
$hotelSheet.Columns("C").Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Swap the tab order so that "review_info" comes before "hotel_info".
$hotelSheet.Move($null, $reviewSheet)
